$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking
# strings (e.g. "1.00", "168.58") are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '64.957.29'
$ws.Range("E2").Value = '  +0.14%  '

$ws.Range("D3").Value = '3.517.21'
$ws.Range("E3").Value = '  -0.23%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").Value = '597.33'
$ws.Range("E5").Value = '  +0.90%  '

$ws.Range("D6").Value = '133.73'
$ws.Range("E6").Value = '  -2.31%  '

$ws.Range("D7").Value = '3.516.44'
$ws.Range("E7").Value = '  -0.32%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").Value = '0.494'
$ws.Range("E9").Value = '  +0.28%  '

$ws.Range("D10").Value = '0.123'
$ws.Range("E10").Value = '  +0.77%  '

$ws.Range("D11").Value = '7.13'
$ws.Range("E11").Value = '  +3.24%  '

$ws.Range("D12").Value = '0.382'
$ws.Range("E12").Value = '  -0.08%  '

$ws.Range("D13").Value = '4.128.43'
$ws.Range("E13").Value = '  +0.11%  '

$ws.Range("D14").Value = '27.33'
$ws.Range("E14").Value = '  +0.83%  '

$ws.Range("D15").Value = '0.0000181'
$ws.Range("E15").Value = '  +0.31%  '

$ws.Range("E16").Value = '  +0.05%  '

$ws.Range("D17").Value = '3.526.60'
$ws.Range("E17").Value = '  +0.11%  '

$ws.Range("D18").Value = '64.958.69'
$ws.Range("E18").Value = '  +0.14%  '

$ws.Range("D19").Value = '9.99'
$ws.Range("E19").Value = '  +0.33%  '

$ws.Range("D20").Value = '14.37'
$ws.Range("E20").Value = '  +1.46%  '

$ws.Range("D21").Value = '5.66'
$ws.Range("E21").Value = '  -2.45%  '

$ws.Range("D22").Value = '390.87'
$ws.Range("E22").Value = '  +0.83%  '

$ws.Range("D23").Value = '0.574'
$ws.Range("E23").Value = '  +0.22%  '

$ws.Range("D24").Value = '3.662.25'
$ws.Range("E24").Value = '  -0.05%  '

$ws.Range("D25").Value = '73.99'
$ws.Range("E25").Value = '  +0.33%  '

$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.01%  '

$ws.Range("D27").Value = '0.0000113'
$ws.Range("E27").Value = '  +0.34%  '

$ws.Range("D28").Value = '1.63'
$ws.Range("E28").Value = '  +21.26%  '

$ws.Range("D29").Value = '7.68'
$ws.Range("E29").Value = '  +1.07%  '

$ws.Range("E30").Value = '  +0.04%  '

$ws.Range("D31").Value = '2.28'
$ws.Range("E31").Value = '  +1.57%  '

$ws.Range("D32").Value = '8.31'
$ws.Range("E32").Value = '  +1.89%  '

$ws.Range("D33").Value = '3.522.94'
$ws.Range("E33").Value = '  -0.38%  '

$ws.Range("D34").Value = '24.04'
$ws.Range("E34").Value = '  +1.42%  '

$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("D36").Value = '0.144'
$ws.Range("E36").Value = '  +1.13%  '

$ws.Range("D37").Value = '1.57'
$ws.Range("E37").Value = '  +1.54%  '

$ws.Range("D38").Value = '5.16'
$ws.Range("E38").Value = '  +4.64%  '

$ws.Range("D39").Value = '168.58'

$ws.Range("D40").Value = '6.81'
$ws.Range("E40").Value = '  -0.17%  '

$ws.Range("D41").Value = '0.0816'
$ws.Range("E41").Value = '  +2.34%  '

$ws.Range("D42").Value = '0.820'
$ws.Range("E42").Value = '  +0.11%  '

$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").Value = '42.69'
$ws.Range("E43").Value = '  +0.50%  '

$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D44").Value = '1.24'
$ws.Range("E44").Value = '  +3.40%  '

$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.23%  '

$ws.Range("D46").Value = '25.26'
$ws.Range("E46").Value = '  -5.19%  '

$ws.Range("D47").Value = '4.41'
$ws.Range("E47").Value = '  +0.10%  '

$ws.Range("D48").Value = '1.65'
$ws.Range("E48").Value = '  -0.39%  '

$ws.Range("D49").Value = '6.88'
$ws.Range("E49").Value = '  +0.53%  '

$ws.Range("D50").Value = '2.405.68'
$ws.Range("E50").Value = '  +0.34%  '

$ws.Range("D51").Value = '0.892'
$ws.Range("E51").Value = '  +5.18%  '
